# Add data for 2023-09-21: update the 2023 (and a few prior-year correction)
# cell values across the Citywide Totals, By Neighborhood summary, and each
# individual neighborhood worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5523
$ws.Range("G3").Value = 8320
$ws.Range("J3").Value = 5890
$ws.Range("I4").Value = 1771
$ws.Range("J4").Value = 1282
$ws.Range("J5").Value = 455
$ws.Range("J6").Value = 7395
$ws.Range("G7").Value = 24696
$ws.Range("I7").Value = 26225
$ws.Range("J7").Value = 20545

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 166
$ws.Range("J6").Value = 148
$ws.Range("G7").Value = 708
$ws.Range("J7").Value = 597
$ws.Range("J8").Value = 1288
$ws.Range("J11").Value = 317
$ws.Range("J12").Value = 41
$ws.Range("J15").Value = 226
$ws.Range("J16").Value = 77
$ws.Range("J19").Value = 600
$ws.Range("J20").Value = 423
$ws.Range("J21").Value = 58
$ws.Range("J27").Value = 123
$ws.Range("J29").Value = 1152
$ws.Range("J33").Value = 943
$ws.Range("J36").Value = 285
$ws.Range("J37").Value = 633
$ws.Range("J41").Value = 127
$ws.Range("J42").Value = 849
$ws.Range("J43").Value = 169
$ws.Range("J47").Value = 158
$ws.Range("J51").Value = 253
$ws.Range("J52").Value = 522
$ws.Range("J54").Value = 399
$ws.Range("J58").Value = 12
$ws.Range("I63").Value = 238
$ws.Range("J63").Value = 77
$ws.Range("J65").Value = 526
$ws.Range("J67").Value = 780
$ws.Range("J68").Value = 38
$ws.Range("J70").Value = 29
$ws.Range("J72").Value = 83
$ws.Range("J76").Value = 297
$ws.Range("J78").Value = 252
$ws.Range("J79").Value = 588
$ws.Range("J83").Value = 419
$ws.Range("J85").Value = 869
$ws.Range("J86").Value = 126
$ws.Range("J88").Value = 222
$ws.Range("J89").Value = 269
$ws.Range("J90").Value = 224
$ws.Range("J93").Value = 95
$ws.Range("J95").Value = 308
$ws.Range("J97").Value = 169
$ws.Range("J98").Value = 144
$ws.Range("J99").Value = 322
$ws.Range("G101").Value = 24696
$ws.Range("I101").Value = 26225
$ws.Range("J101").Value = 20545

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 186
$ws.Range("G3").Value = 268
$ws.Range("J6").Value = 191
$ws.Range("G7").Value = 708
$ws.Range("J7").Value = 597

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 65
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 74
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 228
$ws.Range("J6").Value = 252
$ws.Range("J7").Value = 869

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 160
$ws.Range("J7").Value = 522

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 391
$ws.Range("J4").Value = 74
$ws.Range("J7").Value = 1288

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 125
$ws.Range("J7").Value = 419

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 232
$ws.Range("J3").Value = 306
$ws.Range("J4").Value = 40
$ws.Range("J6").Value = 325
$ws.Range("J7").Value = 943

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 107
$ws.Range("J7").Value = 308

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 215
$ws.Range("J6").Value = 184
$ws.Range("J7").Value = 633

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 155
$ws.Range("J7").Value = 526

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 128
$ws.Range("J7").Value = 322

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 300
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 780

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 99
$ws.Range("J3").Value = 78
$ws.Range("J6").Value = 188
$ws.Range("J7").Value = 399

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 344
$ws.Range("J3").Value = 396
$ws.Range("J6").Value = 304
$ws.Range("J7").Value = 1152

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 151
$ws.Range("J6").Value = 220
$ws.Range("J7").Value = 600

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 189
$ws.Range("J7").Value = 849

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 83
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 166
$ws.Range("J3").Value = 208
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 588

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 147
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 423

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 93
$ws.Range("J7").Value = 285

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 55
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J3").Value = 42
$ws.Range("J4").Value = 12
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J2").Value = 17
$ws.Range("J4").Value = 68
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 65
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 60
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 12
